$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: bump the date serial by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Column D price updates for rows 29-38
$ws.Range("D29").Value = 221.7
$ws.Range("D30").Value = 221.7
$ws.Range("D31").Value = 221.7
$ws.Range("D32").Value = 221.7
$ws.Range("D33").Value = 238.3
$ws.Range("D34").Value = 238.3
$ws.Range("D35").Value = 238.3
$ws.Range("D36").Value = 238.3
$ws.Range("D37").Value = 263
$ws.Range("D38").Value = 263
